$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("Q2").Value = 1.67
$ws.Range("R2").Value = 2.15

# Row 4
$ws.Range("M4").Value = 1.1
$ws.Range("N4").Value = 7
$ws.Range("Q4").Value = 2.5
$ws.Range("R4").Value = 1.5

# Row 5
$ws.Range("Q5").Value = 2.5
$ws.Range("R5").Value = 1.5

# Row 7
$ws.Range("G7").Value = 2.55
$ws.Range("H7").Value = 3.05
$ws.Range("I7").Value = 2.7
$ws.Range("K7").Value = 2.02
$ws.Range("L7").Value = 3.25
$ws.Range("N7").Value = 7.8
$ws.Range("S7").Value = 1.42
$ws.Range("T7").Value = 2.47
$ws.Range("X7").Value = 12
$ws.Range("Z7").Value = 28
$ws.Range("AD7").Value = 5.9
$ws.Range("AF7").Value = 70
$ws.Range("AH7").Value = 8
$ws.Range("AI7").Value = 13.5
$ws.Range("AK7").Value = 32
$ws.Range("AL7").Value = 24
$ws.Range("AN7").Value = 4.4
$ws.Range("AO7").Value = 13.5
$ws.Range("AT7").Value = 2.45
$ws.Range("AW7").Value = 4.55
$ws.Range("AX7").Value = 14.5
$ws.Range("AY7").Value = 22
$ws.Range("AZ7").Value = 65
$ws.Range("BA7").Value = 100
$ws.Range("BB7").Value = 300
